# "Up to date with run 22"
#
# Re-selects a number of cells across sheets (reflecting the user's latest
# navigation), switches the EML milestone-year length used by the
# TimePeriods lookup table from EMLv2 to EMLv3 (which recalculates the
# dependent HLOOKUP/date formulas), and leaves Interpol_Extrapol_Defaults
# as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- Region-Time Slices: just a new selection -------------------------
$wsRegion = $wb.Worksheets.Item("Region-Time Slices")
$wsRegion.Activate()
$wsRegion.Range("D9").Select()

# --- TimePeriods: change the active EML choice + view ------------------
$wsTime = $wb.Worksheets.Item("TimePeriods")
$wsTime.Activate()
$wsTime.Range("B8").Value = "EMLv3"
$wsTime.Range("D12").Select()
$excel.ActiveWindow.Zoom = 122

# --- Defaults: keep selection, but it will lose the active tab ---------
$wsDefaults = $wb.Worksheets.Item("Defaults")
$wsDefaults.Activate()
$wsDefaults.Range("F17").Select()

# --- Comm: new selection -------------------------------------------------
$wsComm = $wb.Worksheets.Item("Comm")
$wsComm.Activate()
$wsComm.Range("K11").Select()

# --- Interpol_Extrapol_Defaults: becomes the active tab, activate last -
$wsInterpol = $wb.Worksheets.Item("Interpol_Extrapol_Defaults")
$wsInterpol.Activate()
$wsInterpol.Range("K19").Select()

Write-Host "Applied run 22 updates"
